{"js": "// Fix the typo \"rusprass\" -> \"surpasses\" in the assessment scale paragraph,\n// and keep the \"_GoBack\" bookmark anchored at the point of the last edit\n// (i.e. move it from its old location, right after \"outlook\", to right\n// after the newly corrected \"surpasses\").\n\nconst misspelled = context.document.body.search(\"rusprass\", { matchCase: true });\nmisspelled.load(\"text\");\nawait context.sync();\n\nif (misspelled.items.length > 0) {\n  const target = misspelled.items[0];\n\n  // Replace the misspelled word in place with the corrected word.\n  target.insertText(\"surpasses\", Word.InsertLocation.replace);\n  await context.sync();\n\n  // Word keeps a single \"_GoBack\" bookmark that marks the location of the\n  // most recent edit. Drop the old one (if present) and re-insert it right\n  // after the text we just fixed.\n  context.document.deleteBookmark(\"_GoBack\");\n  await context.sync();\n\n  const fixed = context.document.body.search(\"surpasses\", { matchCase: true });\n  await context.sync();\n\n  const afterFixed = fixed.items[0].getRange(Word.RangeLocation.after);\n  afterFixed.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Fix the typo \"rusprass\" -> \"surpasses\" in the assessment scale paragraph,\n# and keep the \"_GoBack\" bookmark anchored at the point of the last edit\n# (i.e. move it from its old location, right after \"outlook\", to right\n# after the newly corrected \"surpasses\").\n\n$d = $word.ActiveDocument\n\n# Word tracks the location of the last edit with a single hidden \"_GoBack\"\n# bookmark. Remove the old one (wherever it currently sits) before adding\n# the new one, since bookmark names must stay unique.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Locate and correct the misspelling.\n$target = $d.Content\n$target.Find.Text = \"rusprass\"\nif ($target.Find.Execute()) {\n    $target.Text = \"surpasses\"\n}\n\n# Re-find the corrected word so we can drop a collapsed range right after it,\n# which is where Word re-plants the \"_GoBack\" bookmark after an edit.\n$fixed = $d.Content\n$fixed.Find.Text = \"surpasses\"\nif ($fixed.Find.Execute()) {\n    $insertionPoint = $d.Range($fixed.End, $fixed.End)\n    $d.Bookmarks.Add(\"_GoBack\", $insertionPoint)\n}\n"}
